# Appends 21 new daily rows (2024-08-28 .. 2024-09-25) to the bottom of the
# existing stock-history table on Sheet1, rows 634-654, columns A and C:J.
# (Columns B and I-row-1 style quirks aside, this mirrors the existing
# layout: A=Date(text), C..H=prices/basket value, I=daily return, J=NAV.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date,C,D,E,F,G,H,I,J  (decimal notation only - no scientific notation,
# since the host PowerShell numeric literal parser chokes on the `e` suffix)
$data = @"
2024-08-28,1074.550048828125,700.0999755859375,260.4500122070312,491.7000122070312,1535.849975585938,29293.15045166016,0.0,247.9982131930638
2024-08-29,1121.650024414062,689.1500244140625,253.8000030517578,490.2000122070312,1520.25,29353.35063171387,0.002055094079179154,248.5078728526439
2024-08-30,1111.349975585938,683.3499755859375,256.4500122070312,492.8999938964844,1521.0,29291.19976806641,-0.002117334556700047,247.981698545841
2024-09-02,1092.650024414062,715.0499877929688,251.3500061035156,490.5,1505.25,29322.30020141602,0.00106176713811209,248.2449973642102
2024-09-03,1085.099975585938,710.7999877929688,251.0,488.8500061035156,1509.0,29201.39971923828,-0.004123158188384413,247.2214439706025
2024-09-04,1080.449951171875,722.4000244140625,250.5,484.1499938964844,1488.099975585938,29200.89978027344,-0.00001712037675078922,247.2172114463409
2024-09-05,1069.150024414062,733.8499755859375,251.1499938964844,495.6499938964844,1447.599975585938,29262.3996887207,0.002106096350113556,247.7378747130533
2024-09-06,1049.349975585938,718.9000244140625,247.8000030517578,483.0,1418.050048828125,28702.20024108887,-0.01914400232349252,242.9951802639295
2024-09-09,1038.699951171875,700.1500244140625,243.8999938964844,474.75,1411.849975585938,28242.6496887207,-0.01601098691069303,239.1045876133622
2024-09-10,1035.800048828125,713.4000244140625,248.25,478.7999877929688,1424.449951171875,28522.85046386719,0.009921192885042528,241.4767903467729
2024-09-11,976.2999877929688,725.4000244140625,241.5500030517578,472.2000122070312,1399.599975585938,27922.25028991699,-0.02105680758348599,236.3920600365631
2024-09-12,986.1500244140625,726.0499877929688,246.1499938964844,479.8500061035156,1403.150024414062,28182.30001831055,0.009313351384414074,238.5936623561691
2024-09-13,992.0999755859375,724.25,245.6499938964844,485.3999938964844,1410.949951171875,28258.59951782227,0.002707355306775728,239.2396201741122
2024-09-16,988.4000244140625,733.6500244140625,243.8000030517578,489.9500122070312,1404.550048828125,28313.45072937012,0.00194104493795801,239.7039950278101
2024-09-17,974.9500122070312,745.4000244140625,240.8000030517578,482.2999877929688,1400.25,28196.30033874512,-0.004137623200533361,238.7121902167225
2024-09-18,962.0499877929688,717.5499877929688,235.9499969482422,471.75,1391.300048828125,27572.89979553223,-0.02210930284198537,233.4344301111474
2024-09-19,967.0,728.5,237.5500030517578,459.9500122070312,1374.150024414062,27641.0502166748,0.002471645044516532,234.0113971635512
2024-09-20,970.8499755859375,748.3499755859375,237.8500061035156,466.2999877929688,1380.550048828125,27960.69967651367,0.01156430227264067,236.7175756955934
2024-09-23,971.7999877929688,750.2000122070312,236.4499969482422,471.1499938964844,1375.400024414062,27984.94996643066,0.0008672991090191444,236.9228806380833
2024-09-24,977.2999877929688,735.9000244140625,237.3000030517578,476.7000122070312,1363.699951171875,27912.20024108887,-0.002599601765558408,236.3069754992754
2024-09-25,963.5999755859375,730.0499877929688,238.3500061035156,473.7000122070312,1365.400024414062,27741.79995727539,-0.006104867489544392,234.8643527269973
"@

$lines = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Length

$startRow = 634
$lastRow = $startRow + $rowCount - 1

# --- Column A (dates) -------------------------------------------------
# Writing "2024-08-28" straight into Range.Value lets Excel's input
# parser auto-convert it to a real date serial; the source workbook
# stores these as literal text, so stage each value in a scratch column
# with a leading apostrophe (forces text) and then bring it across with
# PasteSpecial(xlPasteValues), which carries the text over without the
# text-number-format style that a direct text entry would otherwise pick up.
$scratchCol = 30  # column AD, far outside the used range
for ($i = 0; $i -lt $rowCount; $i++) {
    $fields = $lines[$i] -split ","
    $dateText = $fields[0]
    $ws.Cells.Item($startRow + $i, $scratchCol).Value = "'" + $dateText
}
$scratchRange = $ws.Range($ws.Cells.Item($startRow, $scratchCol), $ws.Cells.Item($lastRow, $scratchCol))
$scratchRange.Copy()
$destDates = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($lastRow, 1))
$destDates.PasteSpecial(-4163)  # xlPasteValues
$scratchRange.Clear()

# --- Columns C:J (numbers) --------------------------------------------
$numData = New-Object 'object[,]' $rowCount,8
for ($i = 0; $i -lt $rowCount; $i++) {
    $fields = $lines[$i] -split ","
    for ($c = 0; $c -lt 8; $c++) {
        $numData[$i, $c] = [double]$fields[$c + 1]
    }
}
$destNums = $ws.Range($ws.Cells.Item($startRow, 3), $ws.Cells.Item($lastRow, 10))
$destNums.Value = $numData

$ws.Application.CutCopyMode = $false

Write-Output ("Appended rows " + $startRow + "-" + $lastRow + " to " + $ws.Name)
